$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the new date cells as Text first so Excel stores the literal
# "dd-mm-yyyy" strings (matching the existing column A data) instead of
# auto-converting them into date serial numbers.
$dateRange = $ws.Range("A194:A197")
$dateRange.NumberFormat = "@"

$ws.Cells.Item(194, 1).Value = "02-11-2021"
$ws.Cells.Item(194, 4).Value = 3.28

$ws.Cells.Item(195, 1).Value = "03-11-2021"
$ws.Cells.Item(195, 3).Value = 2.7
$ws.Cells.Item(195, 4).Value = 3.27

$ws.Cells.Item(196, 1).Value = "04-11-2021"
$ws.Cells.Item(196, 4).Value = 3.4

$ws.Cells.Item(197, 1).Value = "05-11-2021"
$ws.Cells.Item(197, 3).Value = 2.68
$ws.Cells.Item(197, 4).Value = 3.34

# Restore the default (unformatted) style on column A so the new rows
# keep the same plain look as the rest of the data (only the header row
# uses a named style).
$dateRange.Style = "Normal"
